$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Daily Project Status" (sheet2) - add new row 23 (Vipin)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Daily Project Status")
$ws2.Range("A23").Value = "Vipin"
$ws2.Range("B23").Value = 9654033937
$ws2.Range("C23").Value = "vipinchauhan247@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("C23"), "mailto:vipinchauhan247@gmail.com")
$ws2.Range("C23").Style = $ws2.Range("C20").Style
$ws2.Range("D23").Value = "Multi purpose Electronic"
$ws2.Range("E23").Value = "18th Feb"
$ws2.Range("F23").Value = 800
$ws2.Range("G23").Value = 300
$ws2.Range("A13").Select()
$ws2.Range("A23").Select()

# ---------------------------------------------------------------------------
# Sheet: "Student Visited In Lab" (sheet3) - add new row 18 (Prakash Pandey)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Student Visited In Lab")
$ws3.Range("A18").Value = "Prakash Pandey"
$ws3.Range("B18").Value = 8130283991
$ws3.Range("C18").Value = "luckyprakash021@gmail.com"
$ws3.Hyperlinks.Add($ws3.Range("C18"), "mailto:luckyprakash021@gmail.com")
$ws3.Range("C18").Style = $ws3.Range("C17").Style
$ws3.Range("D18").Value = "G.N.I.O.T"
$ws3.Range("E18").Value = "Project"
$ws3.Range("D12").Select()

# ---------------------------------------------------------------------------
# Sheet: "Celeab Responsibilities" (sheet4) - only selection moves, no data
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Celeab Responsibilities")
$ws4.Range("A13").Select()

# ---------------------------------------------------------------------------
# Sheet: "Balance" (sheet5) - fix E5:E7 and add new row 8
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Balance")
$ws5.Range("E5").Value = 1049
$ws5.Range("E6").Value = 1049
$ws5.Range("E7").Value = 1049
$ws5.Range("B8").Value = "18th feb 2015"
$ws5.Range("C8").Value = 2300
$ws5.Range("D8").Value = 735
$ws5.Range("E8").Value = 2614
$ws5.Range("F8").Value = 2635
$ws5.Range("F8").Select()

# ---------------------------------------------------------------------------
# Sheet: "Purchase " (sheet6) - add new rows 26 and 27
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Purchase ")
$ws6.Range("A26").Value = "18th feb"
$ws6.Range("B26").Value = "Given to Sweeper"
$ws6.Range("C26").Value = 700
$ws6.Range("A27").Value = "18th feb"
$ws6.Range("B27").Value = "Carbon Paper"
$ws6.Range("C27").Value = 10
$ws6.Range("A13").Select()
$ws6.Range("C27").Select()

# ---------------------------------------------------------------------------
# Sheet: "Rate List" (sheet7) - only selection/view changes, no data
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Rate List")
$ws7.Range("B27").Select()

# ---------------------------------------------------------------------------
# Sheet: "Major Projects Undergoing" (sheet8) - add new row 22 (Sohan)
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Major Projects Undergoing")
$ws8.Range("A22").Value = 22
$ws8.Range("B22").Value = "Variable wind power plant"
$ws8.Range("C22").Value = "18th Feb"
$ws8.Range("D22").Value = "Sohan"
$ws8.Range("E22").Value = "sohanshine@gmail.com"
$ws8.Hyperlinks.Add($ws8.Range("E22"), "mailto:sohanshine@gmail.com")
$ws8.Range("F22").Value = "Sharda"
$ws8.Range("G22").Value = "9136791828"
$ws8.Range("H22").Value = 3000
$ws8.Range("I22").Value = 1000
$ws8.Range("J22").Value = "Within 1 Month"
$ws8.Range("J22").Select()

# ---------------------------------------------------------------------------
# Sheet: "Final Year Projects" (sheet9) - only selection/view changes, no data
# ---------------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item("Final Year Projects")
$ws9.Range("A79").Select()

# ---------------------------------------------------------------------------
# Final active sheet/tab: "Purchase " (tab index 5, 0-based) per workbook.xml
# ---------------------------------------------------------------------------
$ws6.Activate()
$ws6.Range("C27").Select()
